# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# Update OFF sheet (Road row, row 3) with updated target depth data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 217
$wsOff.Range("C3").Value = 152
$wsOff.Range("D3").Value = 63
$wsOff.Range("E3").Value = 28
$wsOff.Range("F3").Value = 5
$wsOff.Range("G3").Value = 6

# Update DEF sheet (Road row, row 3) with updated target depth data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 258
$wsDef.Range("C3").Value = 196
$wsDef.Range("D3").Value = 62
$wsDef.Range("E3").Value = 32
$wsDef.Range("F3").Value = 5
$wsDef.Range("G3").Value = 1
